$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: swap in the "3.0.2.201" driver data
$ws.Range("A3").Value = "Qualcomm Atheros AR9580 Wireless Network Adapter - 3.0.2.201"
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 799
$ws.Range("D3").Value = 92.40000000000001

# Row 4: swap in the "10.1.10.5" driver data
$ws.Range("A4").Value = "Qualcomm Atheros AR9580 Wireless Network Adapter - 10.1.10.5"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 273
$ws.Range("D4").Value = 97.2

# Row 5: update totals (Critical Minutes total changes)
$ws.Range("C5").Value = 1072
